# Add two new trailing columns, I ("I0") and J ("IF"), to the sheet.
# Mirrors the existing H column's header formatting (bold, bordered,
# centered style = style index 1 in the original workbook) by copying
# it across rather than rebuilding it from scratch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy H1's formatting onto I1:J1, then set text ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-81: bulk-write the numeric values for columns I and J ---
$iValues = @(7,8,9,8,7,9,8,8,8,9,8,7,8,9,8,7,6,7,8,8,6,10,7,10,7,8,8,9,7,7,8,8,7,7,8,7,7,8,7,8,7,8,7,8,8,8,6,8,7,9,8,8,8,8,8,8,6,9,7,7,7,7,7,7,7,8,6,8,7,8,8,6,6,8,8,5,5,5,6,3)
$jValues = @(7,8,9,8,7,9,8,9,8,9,9,8,8,9,8,7,7,7,8,8,7,10,7,10,7,8,8,9,7,7,8,8,8,7,8,7,7,8,7,8,7,8,7,8,8,8,7,8,7,9,8,8,8,9,8,9,6,9,8,7,7,7,8,8,7,8,6,8,8,8,8,6,6,8,8,5,5,5,6,3)

$firstRow = 2
$lastRow = 81

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $idx = $r - $firstRow
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]   # column I
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]  # column J
}
